$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.287649333333333
$ws.Range("N2").Value = 3.862948
$ws.Range("O2").Value = 0.1870471291026542
$ws.Range("P2").Value = 0.1870471291026542
$ws.Range("Q2").Value = 28.13484005712089
$ws.Range("R2").Value = 253.213560514088
$ws.Range("S2").Value = 0.009389801215801087
$ws.Range("T2").Value = 0.009389801215801087

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.622908
$ws.Range("N3").Value = 13.868724
$ws.Range("O3").Value = 0.6715350578151914
$ws.Range("P3").Value = 0.6715350578151914
$ws.Range("Q3").Value = 101.0094703672827
$ws.Range("R3").Value = 909.0852333055441
$ws.Range("S3").Value = 0.03371118676120147
$ws.Range("T3").Value = 0.03371118676120147

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9735329999999999
$ws.Range("N4").Value = 2.920599
$ws.Range("O4").Value = 0.1414178130821545
$ws.Range("P4").Value = 0.1414178130821545
$ws.Range("Q4").Value = 21.271470839366
$ws.Range("R4").Value = 191.443237554294
$ws.Range("S4").Value = 0.007099200931792876
$ws.Range("T4").Value = 0.007099200931792876

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.287649333333333
$ws.Range("N5").Value = 3.862948
$ws.Range("O5").Value = 0.1870471291026542
$ws.Range("P5").Value = 0.1870471291026542
$ws.Range("Q5").Value = 495.8125567230699
$ws.Range("R5").Value = 4462.313010507629
$ws.Range("S5").Value = 0.1654738871262716
$ws.Range("T5").Value = 0.1654738871262716

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.622908
$ws.Range("N6").Value = 13.868724
$ws.Range("O6").Value = 0.6715350578151914
$ws.Range("P6").Value = 0.6715350578151914
$ws.Range("Q6").Value = 1780.062145523729
$ws.Range("R6").Value = 16020.55930971357
$ws.Range("S6").Value = 0.5940829826757734
$ws.Range("T6").Value = 0.5940829826757734

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9735329999999999
$ws.Range("N7").Value = 2.920599
$ws.Range("O7").Value = 0.1414178130821545
$ws.Range("P7").Value = 0.1414178130821545
$ws.Range("Q7").Value = 374.861286601021
$ws.Range("R7").Value = 3373.751579409189
$ws.Range("S7").Value = 0.1251072676274963
$ws.Range("T7").Value = 0.1251072676274963

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.287649333333333
$ws.Range("N8").Value = 3.862948
$ws.Range("O8").Value = 0.1870471291026542
$ws.Range("P8").Value = 0.1870471291026542
$ws.Range("Q8").Value = 36.50547538403067
$ws.Range("R8").Value = 328.549278456276
$ws.Range("S8").Value = 0.0121834407605815
$ws.Range("T8").Value = 0.0121834407605815

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.622908
$ws.Range("N9").Value = 13.868724
$ws.Range("O9").Value = 0.6715350578151914
$ws.Range("P9").Value = 0.6715350578151914
$ws.Range("Q9").Value = 131.061656172932
$ws.Range("R9").Value = 1179.554905556388
$ws.Range("S9").Value = 0.04374088837821655
$ws.Range("T9").Value = 0.04374088837821655

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.9735329999999999
$ws.Range("N10").Value = 2.920599
$ws.Range("O10").Value = 0.1414178130821545
$ws.Range("P10").Value = 0.1414178130821545
$ws.Range("Q10").Value = 27.60012687230699
$ws.Range("R10").Value = 248.401141850763
$ws.Range("S10").Value = 0.009211344522865324
$ws.Range("T10").Value = 0.009211344522865324
